$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "Save" column (H) after the existing "sum" column (G), matching
# the header formatting already used by the other header cells in row 1.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Populate the Save values for each data row (2-11).
$saveValues = @(0, 0, 1, 1, 1, 1, 0, 0, 1, 1)
for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}

$excel.CutCopyMode = 0
